$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin names, links) - safe to assign directly
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"

# Numeric-looking text cells (price/volume) - must remain text like original inlineStr cells
# Use an apostrophe prefix to force text entry, then clear the resulting
# number-format / quote-prefix styling so the cell style matches the original (unstyled) cells.
$ws.Range("D2").Value = "'290.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-4.26%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'30.85"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'-6.22%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'4.934"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'-0.28%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.07205"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-8.01%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'1.788"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-10.97%"
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'-2.32%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'3.751"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-1.49%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.8959"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-3.29%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.1666"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'-5.15%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.07718"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-1.09%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.08005"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-7.73%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.03036"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-3.37%"
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = "'-0.23%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.001506"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-0.82%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.04510"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'-1.17%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'0.005723"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'-3.34%"
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'3.480"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'0.40%"
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'2.080"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'-3.46%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'0.3279"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'-0.92%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'0.1299"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'-1.51%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'4.036"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-6.60%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.2101"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'5.51%"
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'-0.88%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.004007"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'-9.96%"
$ws.Range("E25").ClearFormats()
$ws.Range("D39").Value = "'0.01594"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'-8.29%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.04375"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'-8.75%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007347"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-2.05%"
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'-4.02%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.007679"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "'0.002051"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-12.32%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.009209"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-21.43%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00005926"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-5.07%"
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'0.07%"
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = "'173.66%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.003001"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'-3.22%"
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'0.07%"
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'0.07%"
$ws.Range("E51").ClearFormats()
